$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.632.29'
$ws.Range('E2').Value = '  +0.70%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.253.56'
$ws.Range('E3').Value = '  +2.97%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.80'
$ws.Range('E5').Value = '  +1.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.93'
$ws.Range('E6').Value = '  +2.76%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.253.94'
$ws.Range('E8').Value = '  +3.04%  '
$ws.Range('E9').Value = '  +0.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.161'
$ws.Range('E10').Value = '  +2.50%  '
$ws.Range('E11').Value = '  +5.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.504'
$ws.Range('E12').Value = '  -0.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000272'
$ws.Range('E13').Value = '  +2.62%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '39.17'
$ws.Range('E14').Value = '  +2.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.794.53'
$ws.Range('E15').Value = '  +3.18%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.766.72'
$ws.Range('E16').Value = '  +0.82%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.42'
$ws.Range('E17').Value = '  +1.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.258.62'
$ws.Range('E18').Value = '  +3.23%  '
$ws.Range('E19').Value = '  +1.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '509.22'
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.42'
$ws.Range('E21').Value = '  +0.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.751'
$ws.Range('E22').Value = '  +3.67%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.12'
$ws.Range('E23').Value = '  +0.82%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.85'
$ws.Range('E24').Value = '  +1.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.68'
$ws.Range('E25').Value = '  +2.60%  '
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('E27').Value = '  +59.81%  '
$ws.Range('E28').Value = '  +1.90%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.09'
$ws.Range('E29').Value = '  +0.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.41'
$ws.Range('E30').Value = '  +1.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.91'
$ws.Range('E31').Value = '  -4.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.85'
$ws.Range('E32').Value = '  -0.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '28.11'
$ws.Range('E33').Value = '  +1.05%  '
$ws.Range('E34').Value = '  +0.20%  '
$ws.Range('E35').Value = '  -3.46%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.46'
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.40'
$ws.Range('E37').Value = '  +22.96%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0796'
$ws.Range('E38').Value = '  +18.68%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '55.69'
$ws.Range('E39').Value = '  +1.77%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '494.55'
$ws.Range('E40').Value = '  -1.82%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0427'
$ws.Range('E41').Value = '  +2.48%  '
$ws.Range('E42').Value = '  +1.39%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.84'
$ws.Range('E43').Value = '  +0.59%  '
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.51'
$ws.Range('E45').Value = '  +4.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.969.58'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.88'
$ws.Range('E47').Value = '  +4.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.48'
$ws.Range('E48').Value = '  +5.78%  '
$ws.Range('E49').Value = '  +3.08%  '
$ws.Range('E50').Value = '  -0.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '121.60'
$ws.Range('E51').Value = '  +0.37%  '
